$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing values in rows 356-358 (B and D columns) ---
$ws.Cells.Item(356, 2).Value = 2078216858000
$ws.Cells.Item(356, 4).Value = 474869038022.1187

$ws.Cells.Item(357, 2).Value = 2082183969000
$ws.Cells.Item(357, 4).Value = 480662981370.7611

$ws.Cells.Item(358, 2).Value = 2118202312000
$ws.Cells.Item(358, 4).Value = 476707546473.4213

# --- Append new rows 359-361, copying formatting from the row above ---
$ws.Range("A358:D358").Copy()
$ws.Range("A359:D359").PasteSpecial(-4122)
$ws.Range("A358:D358").Copy()
$ws.Range("A360:D360").PasteSpecial(-4122)
$ws.Range("A358:D358").Copy()
$ws.Range("A361:D361").PasteSpecial(-4122)

$ws.Cells.Item(359, 1).Value = 44986
$ws.Cells.Item(359, 2).Value = 2121975670000
$ws.Cells.Item(359, 3).Value = 0.2319647413593134
$ws.Cells.Item(359, 4).Value = 492223537462.3057

$ws.Cells.Item(360, 1).Value = 45017
$ws.Cells.Item(360, 2).Value = 2135028350000
$ws.Cells.Item(360, 3).Value = 0.240610187435336
$ws.Cells.Item(360, 4).Value = 513709571473.2562

$ws.Cells.Item(361, 1).Value = 45047
$ws.Cells.Item(361, 2).Value = 2140971740000
$ws.Cells.Item(361, 3).Value = 0.2359826316783085
$ws.Cells.Item(361, 4).Value = 505232145554.0873
